$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "[핸즈온 머신러닝 2], [GAN 인 액션], [파이썬을 활용한 머신러닝 쿡북], [케라스 창시자에게 배우는 딥러닝], [파이토치로 배우는 자연어 처리] 텐서플로 2.6과 파이토치 1.9 버전 테스트 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/08/24/%ed%95%b8%ec%a6%88%ec%98%a8-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-2-gan-%ec%9d%b8-%ec%95%a1%ec%85%98-%ed%8c%8c%ec%9d%b4%ec%8d%ac%ec%9d%84-%ed%99%9c%ec%9a%a9%ed%95%9c-%eb%a8%b8%ec%8b%a0%eb%9f%ac/"

$ws.Range("D37").Value = "[Paper Review] ERNIE: Enhanced Language Representation with Informative Entities"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1814&mod=document&pageid=1"

$ws.Range("D41").Value = "AI가 게임 개발을 도울 수 있을까"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/ai%eb%8a%94-%ea%b2%8c%ec%9e%84-%ea%b0%9c%eb%b0%9c%ec%9d%84-%eb%8f%84%ec%9a%b8-%ec%88%98-%ec%9e%88%ec%9d%84%ea%b9%8c/"

$ws.Range("D52").Value = "숨은 DS"
